$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update rule_description text for row 30 (shared string update)
$ws.Range("D30").Value = "swh_max_swan > t1 AND anom_swh_p80_waverys > t2"

# 2. Row 2 updates
$ws.Range("F2").Value = 797999.33
$ws.Range("G2").Value = 797999.33
$ws.Range("H2").Value = 1008950
$ws.Range("I2").Value = 1008950

# 3. Row 21 updates
$ws.Range("H21").Value = 4263519.86
$ws.Range("I21").Value = 4263519.86
$ws.Range("J21").Value = 102.52
$ws.Range("K21").Value = 102.52
$ws.Range("L21").Value = 135
$ws.Range("M21").Value = 135

# 4. Row 24 updates
$ws.Range("H24").Value = 322602.86
$ws.Range("I24").Value = 322602.86
$ws.Range("J24").Value = 89.68
$ws.Range("K24").Value = 89.68
$ws.Range("L24").Value = 122
$ws.Range("M24").Value = 122

# 5. Row 27 updates
$ws.Range("H27").Value = 369686.14
$ws.Range("I27").Value = 369686.14
$ws.Range("J27").Value = 121.13
$ws.Range("K27").Value = 121.13
$ws.Range("L27").Value = 158
$ws.Range("M27").Value = 158

# 6. Row 30 updates
$ws.Range("F30").Value = 119993.45
$ws.Range("G30").Value = 119993.45
$ws.Range("H30").Value = 158720
$ws.Range("I30").Value = 158720
$ws.Range("J30").Value = 105.84
$ws.Range("K30").Value = 105.84
$ws.Range("L30").Value = 140
$ws.Range("M30").Value = 140
$ws.Range("R30").Value = 406
$ws.Range("S30").Value = 330
$ws.Range("T30").Value = 1724
$ws.Range("U30").Value = 101
$ws.Range("V30").Value = 0.55
$ws.Range("W30").Value = 0.8
$ws.Range("Y30").Value = 0.65
